# TMTC0032668 test data update - 16 Sep 2024
# "Changed Test Data for LV Activities"
#
# The Contact sheet's sample row (A2:B2) is updated from the generic
# "Test External" / "StandardTestCompany" placeholder values to the new
# LV-Activities-specific values, and that sheet becomes the active sheet
# with A2:B2 selected (mirroring the author's saved UI state).

$wb = $excel.ActiveWorkbook
$contact = $wb.Worksheets.Item("Contact")

$contact.Range("A2").Value = "Activity Test External Contact"
$contact.Range("B2").Value = "ActivityCompany"

# Make Contact the active sheet/selection (was Notification before the edit)
$contact.Activate() | Out-Null
$contact.Range("A2:B2").Select() | Out-Null
